# Populate the "in" column (B) on Sheet1 with the source values that the
# Sheet2 formulas read from, and update the Sheet2 formula for B2 so that
# it too is computed (reading from Sheet1) instead of being a static
# blank/space string.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: add the numeric inputs in column B next to in1/in2/in3
$ws1.Range("B3").Value = 1
$ws1.Range("B4").Value = 5
$ws1.Range("B5").Value = 3

# Sheet2: B2 used to be a literal " " string; now it is a formula that
# reads from Sheet1 (mirrors the existing B3/B4 formulas already there).
$ws2.Range("B2").Formula = "=Sheet1!B3+2"

# Recalculate so the cached <v> formula results are refreshed.
$excel.Calculate()

# Update the selections: Sheet2's selection moves to E3 (no longer the
# active sheet) and Sheet1 becomes the active sheet with D6 selected.
$ws2.Range("E3").Select()
$ws1.Range("D6").Select()
